$wb = $excel.ActiveWorkbook

$src = $wb.Worksheets.Item("tc084")
$src.Copy($null, $src)
$new = $wb.Worksheets.Item("tc084 (2)")
$new.Name = "tc085"

$new.Range("H1").Value = "release"
$new.Range("H2").Value = "Release 06-01-2025"
$new.Columns.Item(7).ColumnWidth = 10.7

$src.Range("A1:G2").Select()

$new.Activate()
$new.Range("A1:H2").Select()
